# team_member.xlsx update:
#  - column E ("职位" / position) for the last four data rows (4-7) is
#    re-labelled: rows 4-5 ("test2","test3") become "部门经理",
#    rows 6-7 ("test4","test5") become "员工" (both were "member").
#  - the sheet's remembered selection moves from F9 to E11.
#  - the workbook window height shrinks from 17260 to 15600.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "部门经理"
$ws.Range("E5").Value = "部门经理"
$ws.Range("E6").Value = "员工"
$ws.Range("E7").Value = "员工"

# Move/update the active selection shown when the sheet is reopened.
$ws.Range("E11").Select()

# Persist the new remembered window height for the workbook view.
$excel.ActiveWindow.Height = 15600
$wb.Windows.Item(1).Height = 15600
